# The deck's single Design/Theme (bound to the Slide Master, and persisted
# as ppt/theme/theme2.xml) currently uses the "Integral" / "Red Violet"
# color scheme. The edit swaps it back to the stock "Office" color scheme
# (the 12 standard theme colors, in clrScheme order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) while leaving the font scheme / format scheme
# untouched, since those are identical between the two themes.
#
# PowerPoint's object model doesn't let automation rename a color scheme or
# theme (that can only be done by picking a Design in the UI / loading a
# .thmx file), so we use the documented, automatable surface for recoloring
# a theme: ThemeColorScheme.Colors(i).RGB.

$p = $ppt.ActivePresentation

$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$tcs = $p.SlideMaster.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeColors.Length; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $packed = $r + ($g * 256) + ($b * 65536)
    $tcs.Colors($i).RGB = $packed
}
